# ezoom_db_objects.xlsx update
#
# A new stored procedure, "relocate_inventory", was documented in the DB
# objects catalog. Its row is inserted right after the existing
# "order_dispatch_display_per_product.sql" entry (row 126) and before the
# "view_ingredient" entry, pushing the three trailing "view" rows down by
# one (old rows 127-129 -> new rows 128-130).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row - shifts rows 127..129 ("view_ingredient",
# "view_lot_in_process", "view_process_step") down to 128..130.
$ws.Rows.Item(127).Insert()

# Populate the new row with the relocate_inventory stored procedure entry.
$ws.Cells.Item(127, 1).Value = "stored procedure"
$ws.Cells.Item(127, 2).Value = "relocate_inventory"
$ws.Cells.Item(127, 3).Value = "relocate items from one location to a new location, depending on destination, either resulting in a merge into existed record (same item, destination, serial no of parent/descendant) or creating a new record"

# Leave the selection where the author ended up after the edit.
$ws.Range("C132").Select() | Out-Null
